$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values per row for columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado), P (Precio $/Kg).
# Values derived from target OOXML state (rows reshuffled / re-dated).
$data = @{
    2  = @{ D = 44592; J = 120; K = 12000; L = 13000; M = 12500; P = 962  }
    3  = @{ D = 44406; J = 160; K = 17000; L = 18000; M = 17500; P = 1346 }
    4  = @{ D = 44320; J = 160; K = 19000; L = 20000; M = 19500; P = 1500 }
    5  = @{ D = 44580; J = 160; K = 11000; L = 12000; M = 11500; P = 885  }
    6  = @{ D = 44159; J = 100; K = 23000; L = 24000; M = 23500; P = 1808 }
    7  = @{ D = 44616; J = 120; K = 19000; L = 20000; M = 19500; P = 1500 }
    8  = @{ D = 44379; J = 120; K = 12000; L = 13000; M = 12667; P = 974  }
    10 = @{ D = 44469; J = 140; K = 13000; L = 14000; M = 13500; P = 1038 }
    11 = @{ D = 44229; J = 120; K = 44000; L = 45000; M = 44500; P = 3423 }
    12 = @{ D = 44397; J = 140; K = 12500; L = 13000; M = 12750; P = 981  }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 4).Value2  = $vals.D   # D - Fecha
    $ws.Cells.Item($row, 10).Value2 = $vals.J   # J - Volumen
    $ws.Cells.Item($row, 11).Value2 = $vals.K   # K - Precio minimo
    $ws.Cells.Item($row, 12).Value2 = $vals.L   # L - Precio maximo
    $ws.Cells.Item($row, 13).Value2 = $vals.M   # M - Precio promedio ponderado
    $ws.Cells.Item($row, 16).Value2 = $vals.P   # P - Precio $/Kg
}
